$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5, pushing current rows 5-7 down to 6-8.
$ws.Rows.Item(5).Insert()

# Copy style (incl. number format) of the date cell from the row above (D6, which
# was D5 before the insert) onto the new D5 so it keeps the date formatting.
$ws.Range("D6").Copy()
$ws.Range("D5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new week's record in row 5.
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(5, 4).Value = 44494
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = 100112001
$ws.Cells.Item(5, 7).Value = "Berenjena"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 120
$ws.Cells.Item(5, 11).Value = 5000
$ws.Cells.Item(5, 12).Value = 6000
$ws.Cells.Item(5, 13).Value = 5500
$ws.Cells.Item(5, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(5, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5, 16).Value = 92
$ws.Cells.Item(5, 17).Value = 60
$ws.Cells.Item(5, 18).Value = "Hortaliza"
